$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-07 Thursday" "2023-09-08 Friday"

Replace-Text "45×22=990" "57×71=4047"
Replace-Text "26×53=1378" "90×62=5580"
Replace-Text "68×30=2040" "18×31=558"
Replace-Text "39×80=3120" "80×76=6080"
Replace-Text "13×92=1196" "75×49=3675"

Replace-Text "40×99=3960" "57×71=4047"
Replace-Text "83×88=7304" "37×17=629"
Replace-Text "76×20=1520" "86×74=6364"
Replace-Text "27×37=999" "49×24=1176"
Replace-Text "41×70=2870" "12×53=636"

Replace-Text "15×22=330" "11×65=715"
Replace-Text "98×39=3822" "76×15=1140"
Replace-Text "88×94=8272" "89×95=8455"
Replace-Text "18×78=1404" "54×23=1242"
Replace-Text "90×42=3780" "41×25=1025"

Replace-Text "19×60=1140" "77×42=3234"
Replace-Text "93×36=3348" "94×50=4700"
Replace-Text "63×56=3528" "54×55=2970"
Replace-Text "75×33=2475" "35×90=3150"
Replace-Text "11×63=693" "11×19=209"

Replace-Text "58×91=5278" "99×62=6138"
Replace-Text "62×81=5022" "45×98=4410"
Replace-Text "42×23=966" "91×48=4368"
Replace-Text "18×43=774" "55×61=3355"
Replace-Text "72×15=1080" "53×19=1007"
